$d = $word.ActiveDocument

# Target: first paragraph ("**ID__AFFARS_5322_topic_7__ID** ") gets a
# paragraph border (5-twip gap on all sides, no line), its left indent
# changes from 120 -> 225 twips, and its text is replaced by the new
# bookmark id, collapsing the two runs (text + trailing space) into one.

$p1 = $d.Paragraphs(1)

# Add the paragraph border (top/left/bottom/right, space=5) and bump the
# left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Range.Borders.DistanceFromTop = 5
$p1.Range.Borders.DistanceFromBottom = 5
$p1.Range.Borders.DistanceFromLeft = 5
$p1.Range.Borders.DistanceFromRight = 5
$p1.Format.LeftIndent = 11.25

# Replace "**ID__AFFARS_5322_topic_7__ID** " (including the trailing
# space run) with "**ID__AFFARS_5322_302__ID**" so the two runs merge
# into a single run with the updated id and no trailing space.
$d.Content.Find.Execute("**ID__AFFARS_5322_topic_7__ID** ", $true, $false,
    $false, $false, $false, $true, 1, $false,
    "**ID__AFFARS_5322_302__ID**", 2) | Out-Null
